$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): values that look like plain numbers must be forced back to
# text after assignment (Excel auto-converts numeric-looking strings), while
# preserving each cell's original (unstyled) formatting.
$priceUpdates = @(
    @{ Addr = "D2"; Value = "66.129.99" },
    @{ Addr = "D3"; Value = "3.316.99" },
    @{ Addr = "D5"; Value = "586.09" },
    @{ Addr = "D6"; Value = "183.06" },
    @{ Addr = "D7"; Value = "0.649" },
    @{ Addr = "D9"; Value = "3.316.64" },
    @{ Addr = "D13"; Value = "3.895.07" },
    @{ Addr = "D15"; Value = "66.182.45" },
    @{ Addr = "D16"; Value = "26.14" },
    @{ Addr = "D18"; Value = "3.299.83" },
    @{ Addr = "D19"; Value = "424.96" },
    @{ Addr = "D20"; Value = "5.55" },
    @{ Addr = "D21"; Value = "13.17" },
    @{ Addr = "D23"; Value = "71.74" },
    @{ Addr = "D26"; Value = "3.463.59" },
    @{ Addr = "D33"; Value = "22.40" },
    @{ Addr = "D38"; Value = "160.86" },
    @{ Addr = "D39"; Value = "1.43" },
    @{ Addr = "D40"; Value = "2.883.24" },
    @{ Addr = "D42"; Value = "26.37" },
    @{ Addr = "D44"; Value = "4.31" },
    @{ Addr = "D45"; Value = "39.95" },
    @{ Addr = "D46"; Value = "0.0662" },
    @{ Addr = "D47"; Value = "5.88" },
    @{ Addr = "D48"; Value = "2.29" },
    @{ Addr = "D49"; Value = "23.17" },
    @{ Addr = "D50"; Value = "313.82" }
)

foreach ($u in $priceUpdates) {
    $cell = $ws.Range($u.Addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = $origStyle
}

# Volume(1h) column (E): percentage strings already contain spaces/%, so they
# never get reinterpreted as numbers - plain assignment is safe.
$volumeUpdates = @(
    @{ Addr = "E2"; Value = "  -0.66%  " },
    @{ Addr = "E3"; Value = "  -0.49%  " },
    @{ Addr = "E4"; Value = "  +0.01%  " },
    @{ Addr = "E5"; Value = "  +2.20%  " },
    @{ Addr = "E6"; Value = "  +1.24%  " },
    @{ Addr = "E7"; Value = "  +4.70%  " },
    @{ Addr = "E8"; Value = "  +0.04%  " },
    @{ Addr = "E9"; Value = "  -0.39%  " },
    @{ Addr = "E10"; Value = "  -3.14%  " },
    @{ Addr = "E11"; Value = "  +2.38%  " },
    @{ Addr = "E12"; Value = "  -0.30%  " },
    @{ Addr = "E13"; Value = "  -0.40%  " },
    @{ Addr = "E14"; Value = "  -2.87%  " },
    @{ Addr = "E15"; Value = "  -0.72%  " },
    @{ Addr = "E16"; Value = "  -3.08%  " },
    @{ Addr = "E17"; Value = "  -2.46%  " },
    @{ Addr = "E18"; Value = "  -1.43%  " },
    @{ Addr = "E19"; Value = "  -2.90%  " },
    @{ Addr = "E20"; Value = "  -2.29%  " },
    @{ Addr = "E21"; Value = "  -2.65%  " },
    @{ Addr = "E22"; Value = "  -2.63%  " },
    @{ Addr = "E23"; Value = "  -2.40%  " },
    @{ Addr = "E24"; Value = "  +0.00%  " },
    @{ Addr = "E25"; Value = "  +0.47%  " },
    @{ Addr = "E26"; Value = "  -0.72%  " },
    @{ Addr = "E27"; Value = "  -0.90%  " },
    @{ Addr = "E28"; Value = "  +4.67%  " },
    @{ Addr = "E29"; Value = "  -3.68%  " },
    @{ Addr = "E30"; Value = "  -1.18%  " },
    @{ Addr = "E31"; Value = "  +0.10%  " },
    @{ Addr = "E32"; Value = "  -2.21%  " },
    @{ Addr = "E33"; Value = "  -1.68%  " },
    @{ Addr = "E34"; Value = "  +0.08%  " },
    @{ Addr = "E35"; Value = "  -2.45%  " },
    @{ Addr = "E36"; Value = "  -3.27%  " },
    @{ Addr = "E37"; Value = "  -4.20%  " },
    @{ Addr = "E38"; Value = "  -1.08%  " },
    @{ Addr = "E39"; Value = "  -3.49%  " },
    @{ Addr = "E40"; Value = "  +2.41%  " },
    @{ Addr = "E41"; Value = "  -2.05%  " },
    @{ Addr = "E42"; Value = "  -3.80%  " },
    @{ Addr = "E43"; Value = "  -3.60%  " },
    @{ Addr = "E44"; Value = "  -2.50%  " },
    @{ Addr = "E45"; Value = "  -0.55%  " },
    @{ Addr = "E46"; Value = "  -0.55%  " },
    @{ Addr = "E47"; Value = "  -5.15%  " },
    @{ Addr = "E48"; Value = "  -2.86%  " },
    @{ Addr = "E49"; Value = "  -5.04%  " },
    @{ Addr = "E50"; Value = "  -2.04%  " },
    @{ Addr = "E51"; Value = "  -0.63%  " }
)

foreach ($u in $volumeUpdates) {
    $ws.Range($u.Addr).Value = $u.Value
}
